$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(32, 8).Value = 5166.3335  # H32 was 5142.4287
$ws.Cells.Item(32, 10).Value = 7499.5  # J32 was 6666
$ws.Cells.Item(32, 12).Value = 7499.5  # L32 was 6666
$ws.Cells.Item(32, 14).Value = -8151.5  # N32 was -7318
$ws.Cells.Item(113, 8).Value = 9972.625  # H113 was 8617.817999999999
$ws.Cells.Item(113, 9).Value = 16249.667  # I113 was 11152.8
$ws.Cells.Item(113, 10).Value = 6206.4  # J113 was 6505.3335
$ws.Cells.Item(113, 11).Value = 16249.667  # K113 was 11152.8
$ws.Cells.Item(113, 12).Value = 6206.4  # L113 was 6505.3335
$ws.Cells.Item(113, 13).Value = -12995.667  # M113 was -7898.799999999999
$ws.Cells.Item(113, 14).Value = -12714.4  # N113 was -13013.3335
$ws.Cells.Item(134, 8).Value = 89082.5  # H134 was 89065.8
$ws.Cells.Item(134, 10).Value = 89082.5  # J134 was 89065.8
$ws.Cells.Item(134, 12).Value = 89082.5  # L134 was 89065.8
$ws.Cells.Item(134, 14).Value = -99222.5  # N134 was -99205.8
$ws.Cells.Item(137, 8).Value = 7575.6587  # H137 was 7740.2
$ws.Cells.Item(137, 9).Value = 10742.292  # I137 was 11166.131
$ws.Cells.Item(137, 11).Value = 32226.876  # K137 was 33498.393
$ws.Cells.Item(137, 13).Value = -29676.876  # M137 was -30948.393
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 3798.9  # H2 was 3079.3572
$ws.Cells.Item(2, 9).Value = 3498.1667  # I2 was 2611.1
$ws.Cells.Item(2, 11).Value = 3498.1667  # K2 was 2611.1
$ws.Cells.Item(2, 13).Value = -3385.1667  # M2 was -2498.1
$ws.Cells.Item(19, 8).Value = 20008  # H19 was 0
$ws.Cells.Item(19, 9).Value = 20008  # I19 was 0
$ws.Cells.Item(19, 11).Value = 20008  # K19 was 0
$ws.Cells.Item(19, 13).Value = -19779  # M19 was None
$ws.Cells.Item(45, 8).Value = 5833  # H45 was 6557.3076
$ws.Cells.Item(45, 9).Value = 4916.222  # I45 was 5437
$ws.Cells.Item(45, 10).Value = 7208.1665  # J45 was 8349.799999999999
$ws.Cells.Item(45, 11).Value = 4916.222  # K45 was 5437
$ws.Cells.Item(45, 12).Value = 7208.1665  # L45 was 8349.799999999999
$ws.Cells.Item(45, 13).Value = -4539.222  # M45 was -5060
$ws.Cells.Item(45, 14).Value = -7962.1665  # N45 was -9103.799999999999
$ws.Cells.Item(74, 8).Value = 5997.5713  # H74 was 6241.4736
$ws.Cells.Item(74, 9).Value = 3880.923  # I74 was 3919.1
$ws.Cells.Item(74, 10).Value = 9437.125  # J74 was 8821.888999999999
$ws.Cells.Item(74, 11).Value = 3880.923  # K74 was 3919.1
$ws.Cells.Item(74, 12).Value = 9437.125  # L74 was 8821.888999999999
$ws.Cells.Item(74, 13).Value = -3006.923  # M74 was -3045.1
$ws.Cells.Item(74, 14).Value = -11185.125  # N74 was -10569.889
$ws.Cells.Item(77, 8).Value = 5997.5713  # H77 was 6241.4736
$ws.Cells.Item(77, 9).Value = 3880.923  # I77 was 3919.1
$ws.Cells.Item(77, 10).Value = 9437.125  # J77 was 8821.888999999999
$ws.Cells.Item(77, 11).Value = 19404.615  # K77 was 19595.5
$ws.Cells.Item(77, 12).Value = 47185.625  # L77 was 44109.44499999999
$ws.Cells.Item(77, 13).Value = -15036.615  # M77 was -15227.5
$ws.Cells.Item(77, 14).Value = -55921.625  # N77 was -52845.44499999999
$ws.Cells.Item(97, 8).Value = 39250.668  # H97 was 16157.739
$ws.Cells.Item(97, 9).Value = 27317  # I97 was 8538.929
$ws.Cells.Item(97, 10).Value = 48797.6  # J97 was 28009.223
$ws.Cells.Item(97, 11).Value = 27317  # K97 was 8538.929
$ws.Cells.Item(97, 12).Value = 48797.6  # L97 was 28009.223
$ws.Cells.Item(97, 13).Value = -26821  # M97 was -8042.929
$ws.Cells.Item(97, 14).Value = -49789.6  # N97 was -29001.223
$ws.Cells.Item(110, 8).Value = 2395.5557  # H110 was 2143.7693
$ws.Cells.Item(110, 9).Value = 1651.4286  # I110 was 1624.4546
$ws.Cells.Item(110, 11).Value = 1651.4286  # K110 was 1624.4546
$ws.Cells.Item(110, 13).Value = 393.5714  # M110 was 420.5454
$ws.Cells.Item(116, 8).Value = 3798.9  # H116 was 3079.3572
$ws.Cells.Item(116, 9).Value = 3498.1667  # I116 was 2611.1
$ws.Cells.Item(116, 11).Value = 3498.1667  # K116 was 2611.1
$ws.Cells.Item(116, 13).Value = -1204.1667  # M116 was -317.0999999999999
$ws.Cells.Item(122, 9).Value = 3789  # I122 was 3789.05
$ws.Cells.Item(122, 11).Value = 11367  # K122 was 11367.15
$ws.Cells.Item(122, 13).Value = -8917  # M122 was -8917.150000000001
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 3798.9  # H3 was 3079.3572
$ws.Cells.Item(3, 9).Value = 3498.1667  # I3 was 2611.1
$ws.Cells.Item(3, 11).Value = 3498.1667  # K3 was 2611.1
$ws.Cells.Item(3, 13).Value = -3384.1667  # M3 was -2497.1
$ws.Cells.Item(22, 8).Value = 216.5  # H22 was 184.85715
$ws.Cells.Item(22, 9).Value = 237.57143  # I22 was 204.16667
$ws.Cells.Item(22, 11).Value = 237.57143  # K22 was 204.16667
$ws.Cells.Item(22, 13).Value = -64.57142999999999  # M22 was -31.16667000000001
$ws.Cells.Item(86, 8).Value = 8201.947  # H86 was 8186.1577
$ws.Cells.Item(86, 9).Value = 5704.4287  # I86 was 5683
$ws.Cells.Item(86, 11).Value = 5704.4287  # K86 was 5683
$ws.Cells.Item(86, 13).Value = -4581.4287  # M86 was -4560
$ws.Cells.Item(89, 8).Value = 8201.947  # H89 was 8186.1577
$ws.Cells.Item(89, 9).Value = 5704.4287  # I89 was 5683
$ws.Cells.Item(89, 11).Value = 28522.1435  # K89 was 28415
$ws.Cells.Item(89, 13).Value = -22906.1435  # M89 was -22799
$ws.Cells.Item(94, 8).Value = 3034.889  # H94 was 3209.0588
$ws.Cells.Item(94, 9).Value = 2158.7856  # I94 was 2319.1538
$ws.Cells.Item(94, 11).Value = 2158.7856  # K94 was 2319.1538
$ws.Cells.Item(94, 13).Value = -1707.7856  # M94 was -1868.1538
$ws.Cells.Item(99, 8).Value = 13965.088  # H99 was 15139.161
$ws.Cells.Item(99, 9).Value = 18152.738  # I99 was 20600.7
$ws.Cells.Item(99, 11).Value = 18152.738  # K99 was 20600.7
$ws.Cells.Item(99, 13).Value = -16654.738  # M99 was -19102.7
$ws.Cells.Item(105, 8).Value = 6448.625  # H105 was 5135.846
$ws.Cells.Item(105, 9).Value = 4098.1665  # I105 was 3615.0908
$ws.Cells.Item(105, 11).Value = 4098.1665  # K105 was 3615.0908
$ws.Cells.Item(105, 13).Value = -2351.1665  # M105 was -1868.0908
$ws.Cells.Item(117, 8).Value = 0  # H117 was 50000
$ws.Cells.Item(117, 10).Value = 0  # J117 was 50000
$ws.Cells.Item(117, 12).Value = 0  # L117 was 50000
$ws.Cells.Item(117, 14).ClearContents()  # N117 was -59178
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 1180.3572  # H7 was 1103
$ws.Cells.Item(7, 9).Value = 1556.5  # I7 was 1416.8182
$ws.Cells.Item(7, 11).Value = 1556.5  # K7 was 1416.8182
$ws.Cells.Item(7, 13).Value = -1443.5  # M7 was -1303.8182
$ws.Cells.Item(14, 8).Value = 225015760  # H14 was 180013600
$ws.Cells.Item(14, 10).Value = 450004500  # J14 was 300004670
$ws.Cells.Item(14, 12).Value = 450004500  # L14 was 300004670
$ws.Cells.Item(14, 14).Value = -450004840  # N14 was -300005010
$ws.Cells.Item(16, 8).Value = 2019  # H16 was 1826.7778
$ws.Cells.Item(16, 9).Value = 2103.3333  # I16 was 1830.3334
$ws.Cells.Item(16, 10).Value = 1892.5  # J16 was 1819.6666
$ws.Cells.Item(16, 11).Value = 2103.3333  # K16 was 1830.3334
$ws.Cells.Item(16, 12).Value = 1892.5  # L16 was 1819.6666
$ws.Cells.Item(16, 13).Value = -1816.3333  # M16 was -1543.3334
$ws.Cells.Item(16, 14).Value = -2466.5  # N16 was -2393.6666
$ws.Cells.Item(31, 8).Value = 4096.1875  # H31 was 4353.077
$ws.Cells.Item(31, 9).Value = 3076.6667  # I31 was 3391.5715
$ws.Cells.Item(31, 10).Value = 5407  # J31 was 5474.8335
$ws.Cells.Item(31, 11).Value = 3076.6667  # K31 was 3391.5715
$ws.Cells.Item(31, 12).Value = 5407  # L31 was 5474.8335
$ws.Cells.Item(31, 13).Value = -2781.6667  # M31 was -3096.5715
$ws.Cells.Item(31, 14).Value = -5997  # N31 was -6064.8335
$ws.Cells.Item(34, 8).Value = 4096.1875  # H34 was 4353.077
$ws.Cells.Item(34, 9).Value = 3076.6667  # I34 was 3391.5715
$ws.Cells.Item(34, 10).Value = 5407  # J34 was 5474.8335
$ws.Cells.Item(34, 11).Value = 3076.6667  # K34 was 3391.5715
$ws.Cells.Item(34, 12).Value = 5407  # L34 was 5474.8335
$ws.Cells.Item(34, 13).Value = -2874.6667  # M34 was -3189.5715
$ws.Cells.Item(34, 14).Value = -5811  # N34 was -5878.8335
$ws.Cells.Item(113, 8).Value = 2019  # H113 was 1826.7778
$ws.Cells.Item(113, 9).Value = 2103.3333  # I113 was 1830.3334
$ws.Cells.Item(113, 10).Value = 1892.5  # J113 was 1819.6666
$ws.Cells.Item(113, 11).Value = 2103.3333  # K113 was 1830.3334
$ws.Cells.Item(113, 12).Value = 1892.5  # L113 was 1819.6666
$ws.Cells.Item(113, 13).Value = 66.66670000000022  # M113 was 339.6666
$ws.Cells.Item(113, 14).Value = -6232.5  # N113 was -6159.6666
$ws.Cells.Item(132, 8).Value = 25028.889  # H132 was 16317.179
$ws.Cells.Item(132, 9).Value = 1977  # I132 was 1232.0555
$ws.Cells.Item(132, 11).Value = 5931  # K132 was 3696.1665
$ws.Cells.Item(132, 13).Value = -3401  # M132 was -1166.1665
$ws.Cells.Item(134, 8).Value = 5746.8  # H134 was 5965.579
$ws.Cells.Item(134, 9).Value = 3146.8462  # I134 was 3276.5833
$ws.Cells.Item(134, 11).Value = 9440.5386  # K134 was 9829.749899999999
$ws.Cells.Item(134, 13).Value = -6905.5386  # M134 was -7294.749899999999
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(97, 8).Value = 61049.285  # H97 was 45589.473
$ws.Cells.Item(97, 9).Value = 86671.28999999999  # I97 was 47784.54
$ws.Cells.Item(97, 10).Value = 35427.285  # J97 was 40833.5
$ws.Cells.Item(97, 11).Value = 260013.87  # K97 was 143353.62
$ws.Cells.Item(97, 12).Value = 106281.855  # L97 was 122500.5
$ws.Cells.Item(97, 13).Value = -259517.87  # M97 was -142857.62
$ws.Cells.Item(97, 14).Value = -107273.855  # N97 was -123492.5
$ws.Cells.Item(132, 8).Value = 64387.625  # H132 was 43575
$ws.Cells.Item(132, 10).Value = 169999.67  # J132 was 73971.14
$ws.Cells.Item(132, 12).Value = 1529997.03  # L132 was 665740.26
$ws.Cells.Item(132, 14).Value = -1535057.03  # N132 was -670800.26
$ws.Cells.Item(139, 8).Value = 2730835.5  # H139 was 1878715
$ws.Cells.Item(139, 9).Value = 3003419  # I139 was 2730449
$ws.Cells.Item(139, 10).Value = 5000  # J139 was 4900
$ws.Cells.Item(139, 11).Value = 9010257  # K139 was 8191347
$ws.Cells.Item(139, 12).Value = 15000  # L139 was 14700
$ws.Cells.Item(139, 13).Value = -9005117  # M139 was -8186207
$ws.Cells.Item(139, 14).Value = -25280  # N139 was -24980
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 1271.3572  # H2 was 1206.6
$ws.Cells.Item(2, 9).Value = 1553.091  # I2 was 1448.6666
$ws.Cells.Item(2, 11).Value = 1553.091  # K2 was 1448.6666
$ws.Cells.Item(2, 13).Value = -1440.091  # M2 was -1335.6666
$ws.Cells.Item(10, 8).Value = 14500  # H10 was 0
$ws.Cells.Item(10, 9).Value = 14500  # I10 was 0
$ws.Cells.Item(10, 11).Value = 14500  # K10 was 0
$ws.Cells.Item(10, 13).Value = -14331  # M10 was None
$ws.Cells.Item(80, 8).Value = 9524.125  # H80 was 9624
$ws.Cells.Item(80, 9).Value = 15000  # I80 was 19000
$ws.Cells.Item(80, 10).Value = 4048.25  # J80 was 3998.4
$ws.Cells.Item(80, 11).Value = 15000  # K80 was 19000
$ws.Cells.Item(80, 12).Value = 4048.25  # L80 was 3998.4
$ws.Cells.Item(80, 13).Value = -14002  # M80 was -18002
$ws.Cells.Item(80, 14).Value = -6044.25  # N80 was -5994.4
$ws.Cells.Item(83, 8).Value = 9524.125  # H83 was 9624
$ws.Cells.Item(83, 9).Value = 15000  # I83 was 19000
$ws.Cells.Item(83, 10).Value = 4048.25  # J83 was 3998.4
$ws.Cells.Item(83, 11).Value = 75000  # K83 was 95000
$ws.Cells.Item(83, 12).Value = 20241.25  # L83 was 19992
$ws.Cells.Item(83, 13).Value = -70008  # M83 was -90008
$ws.Cells.Item(83, 14).Value = -30225.25  # N83 was -29976
$ws.Cells.Item(102, 8).Value = 9910.166999999999  # H102 was 10363.706
$ws.Cells.Item(102, 9).Value = 13115.25  # I102 was 14107.546
$ws.Cells.Item(102, 11).Value = 13115.25  # K102 was 14107.546
$ws.Cells.Item(102, 13).Value = -11493.25  # M102 was -12485.546
$ws.Cells.Item(132, 8).Value = 3656.8948  # H132 was 2828.4614
$ws.Cells.Item(132, 9).Value = 2981.3076  # I132 was 2369.3684
$ws.Cells.Item(132, 10).Value = 5120.6665  # J132 was 4074.5715
$ws.Cells.Item(132, 11).Value = 8943.9228  # K132 was 7108.1052
$ws.Cells.Item(132, 12).Value = 15361.9995  # L132 was 12223.7145
$ws.Cells.Item(132, 13).Value = -6413.9228  # M132 was -4578.1052
$ws.Cells.Item(132, 14).Value = -20421.9995  # N132 was -17283.7145
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(11, 8).Value = 0  # H11 was 1300
$ws.Cells.Item(11, 10).Value = 0  # J11 was 1300
$ws.Cells.Item(11, 12).Value = 0  # L11 was 1300
$ws.Cells.Item(11, 14).ClearContents()  # N11 was -1580
$ws.Cells.Item(16, 8).Value = 10562.6  # H16 was 9561.058999999999
$ws.Cells.Item(16, 9).Value = 10562.6  # I16 was 9561.058999999999
$ws.Cells.Item(16, 11).Value = 10562.6  # K16 was 9561.058999999999
$ws.Cells.Item(16, 13).Value = -10392.6  # M16 was -9391.058999999999
$ws.Cells.Item(24, 8).Value = 15364  # H24 was 21501.5
$ws.Cells.Item(24, 9).Value = 14833.833  # I24 was 19402
$ws.Cells.Item(24, 10).Value = 17749.75  # J24 was 31999
$ws.Cells.Item(24, 11).Value = 14833.833  # K24 was 19402
$ws.Cells.Item(24, 12).Value = 17749.75  # L24 was 31999
$ws.Cells.Item(24, 13).Value = -14490.833  # M24 was -19059
$ws.Cells.Item(24, 14).Value = -18435.75  # N24 was -32685
$ws.Cells.Item(55, 8).Value = 1226.7222  # H55 was 1368.5
$ws.Cells.Item(55, 9).Value = 387.54544  # I55 was 424.8
$ws.Cells.Item(55, 10).Value = 2545.4285  # J55 was 2941.3333
$ws.Cells.Item(55, 11).Value = 387.54544  # K55 was 424.8
$ws.Cells.Item(55, 12).Value = 2545.4285  # L55 was 2941.3333
$ws.Cells.Item(55, 13).Value = -214.54544  # M55 was -251.8
$ws.Cells.Item(55, 14).Value = -2891.4285  # N55 was -3287.3333
$ws.Cells.Item(61, 8).Value = 18237.75  # H61 was 15191
$ws.Cells.Item(61, 9).Value = 15317  # I61 was 12238.75
$ws.Cells.Item(61, 11).Value = 15317  # K61 was 12238.75
$ws.Cells.Item(61, 13).Value = -15115  # M61 was -12036.75
$ws.Cells.Item(113, 8).Value = 18237.75  # H113 was 15191
$ws.Cells.Item(113, 9).Value = 15317  # I113 was 12238.75
$ws.Cells.Item(113, 11).Value = 15317  # K113 was 12238.75
$ws.Cells.Item(113, 13).Value = -13147  # M113 was -10068.75
$ws.Cells.Item(122, 8).Value = 5928.375  # H122 was 5691.4
$ws.Cells.Item(122, 9).Value = 5928.375  # I122 was 5691.4
$ws.Cells.Item(122, 11).Value = 17785.125  # K122 was 17074.2
$ws.Cells.Item(122, 13).Value = -15335.125  # M122 was -14624.2
$ws.Cells.Item(136, 8).Value = 8959.733  # H136 was 7777.1113
$ws.Cells.Item(136, 9).Value = 0  # I136 was 598
$ws.Cells.Item(136, 10).Value = 8959.733  # J136 was 8199.412
$ws.Cells.Item(136, 11).Value = 0  # K136 was 1794
$ws.Cells.Item(136, 12).Value = 26879.199  # L136 was 24598.236
$ws.Cells.Item(136, 13).ClearContents()  # M136 was 756
$ws.Cells.Item(136, 14).Value = -31979.199  # N136 was -29698.236
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(107, 8).Value = 23086.357  # H107 was 24832.615
$ws.Cells.Item(107, 10).Value = 75297.5  # J107 was 100268.336
$ws.Cells.Item(107, 12).Value = 225892.5  # L107 was 300805.008
$ws.Cells.Item(107, 14).Value = -229732.5  # N107 was -304645.008
$ws.Cells.Item(113, 8).Value = 2310.25  # H113 was 2360.375
$ws.Cells.Item(113, 9).Value = 1074.7693  # I113 was 1136.4615
$ws.Cells.Item(113, 11).Value = 3224.3079  # K113 was 3409.3845
$ws.Cells.Item(113, 13).Value = -1054.3079  # M113 was -1239.3845
$ws.Cells.Item(136, 8).Value = 3154.6667  # H136 was 2000.0605
$ws.Cells.Item(136, 9).Value = 2178.2  # I136 was 1240
$ws.Cells.Item(136, 11).Value = 6534.599999999999  # K136 was 3720
$ws.Cells.Item(136, 13).Value = -3984.599999999999  # M136 was -1170
